$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: move "Scalpel Accuracy:" label from C3 to E3, and the value from D3 to F3
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = $null
$ws.Range("E3").Value = "Scalpel Accuracy:"
$ws.Range("F3").Value = 100

# Row 4: fix label text in E4
$ws.Range("E4").Value = "Accuracy vs PyType"
